# Generate Report for Handback
#
# The handback transform for the c8924308-... file failed because the
# handback file name did not match the handoff file name. Update the
# status on the Overview/zh-cn/de-de sheets and record the error detail
# for each locale.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

$zhError = "Handback file name: vc523uqh.5sm is different with handoff file name: c8924308-46b6-4f35-9200-1dd3ff43f2b2.6fbce726e5326730a2ca5801e873972ec533abc5.zh-cn."
$deError = "Handback file name: vc523uqh.5sm is different with handoff file name: c8924308-46b6-4f35-9200-1dd3ff43f2b2.6fbce726e5326730a2ca5801e873972ec533abc5.de-de."

# Update the status for the c8924308 file (row 3) on every sheet that
# reports it: Overview (zh-cn column B, de-de column C) and the two
# per-locale detail sheets (column C).
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Record the handback/handoff file name mismatch in the "Error Detail"
# column (K) for row 3 on each locale sheet.
$zhcn.Range("K3").Value = $zhError
$dede.Range("K3").Value = $deError
